$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D store plain numeric-looking text (e.g. "13.50", "1.98")
# that must remain exact text, not be reinterpreted as numbers. Force text
# format before assignment, then restore the default "Normal" style so the
# cell keeps no explicit style index (matching the original formatting).
$dCells = @(
    "D2"
    "D3"
    "D5"
    "D6"
    "D8"
    "D9"
    "D10"
    "D11"
    "D12"
    "D13"
    "D14"
    "D15"
    "D16"
    "D17"
    "D18"
    "D19"
    "D20"
    "D21"
    "D22"
    "D23"
    "D24"
    "D25"
    "D26"
    "D27"
    "D28"
    "D33"
    "D35"
    "D37"
    "D39"
    "D40"
    "D41"
    "D43"
    "D44"
    "D45"
    "D46"
    "D47"
    "D48"
    "D50"
)
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '60.666.83'
$ws.Range("E2").Value = '  -3.71%  '
$ws.Range("D3").Value = '2.899.39'
$ws.Range("E3").Value = '  -4.25%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '588.78'
$ws.Range("E5").Value = '  -0.96%  '
$ws.Range("D6").Value = '144.12'
$ws.Range("E6").Value = '  -6.05%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '0.503'
$ws.Range("E8").Value = '  -1.77%  '
$ws.Range("D9").Value = '2.898.76'
$ws.Range("E9").Value = '  -4.14%  '
$ws.Range("D10").Value = '6.68'
$ws.Range("E10").Value = '  -5.99%  '
$ws.Range("D11").Value = '0.143'
$ws.Range("E11").Value = '  -4.47%  '
$ws.Range("D12").Value = '0.443'
$ws.Range("E12").Value = '  -4.32%  '
$ws.Range("D13").Value = '0.0000225'
$ws.Range("E13").Value = '  -3.57%  '
$ws.Range("D14").Value = '33.32'
$ws.Range("E14").Value = '  -6.50%  '
$ws.Range("D15").Value = '0.127'
$ws.Range("E15").Value = '  +1.51%  '
$ws.Range("D16").Value = '3.380.60'
$ws.Range("E16").Value = '  -4.22%  '
$ws.Range("D17").Value = '60.646.82'
$ws.Range("E17").Value = '  -3.86%  '
$ws.Range("D18").Value = '6.70'
$ws.Range("E18").Value = '  -5.22%  '
$ws.Range("D19").Value = '2.903.65'
$ws.Range("E19").Value = '  -4.18%  '
$ws.Range("D20").Value = '428.47'
$ws.Range("E20").Value = '  -4.51%  '
$ws.Range("D21").Value = '13.50'
$ws.Range("E21").Value = '  -5.16%  '
$ws.Range("D22").Value = '0.681'
$ws.Range("E22").Value = '  -1.90%  '
$ws.Range("D23").Value = '7.05'
$ws.Range("E23").Value = '  -6.13%  '
$ws.Range("D24").Value = '81.98'
$ws.Range("E24").Value = '  -1.44%  '
$ws.Range("D25").Value = '10.79'
$ws.Range("E25").Value = '  -6.28%  '
$ws.Range("D26").Value = '2.21'
$ws.Range("E26").Value = '  -5.53%  '
$ws.Range("D27").Value = '11.93'
$ws.Range("E27").Value = '  -3.46%  '
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("E29").Value = '  -2.21%  '
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("E31").Value = '  -3.53%  '
$ws.Range("E32").Value = '  -7.17%  '
$ws.Range("D33").Value = '26.46'
$ws.Range("E33").Value = '  -4.25%  '
$ws.Range("E34").Value = '  -3.75%  '
$ws.Range("D35").Value = '0.0₃0846'
$ws.Range("E35").Value = '  -3.00%  '
$ws.Range("E36").Value = '  -3.62%  '
$ws.Range("D37").Value = '5.58'
$ws.Range("E37").Value = '  -5.19%  '
$ws.Range("E38").Value = '  -4.88%  '
$ws.Range("D39").Value = '49.41'
$ws.Range("E39").Value = '  -2.19%  '
$ws.Range("D40").Value = '0.124'
$ws.Range("E40").Value = '  -4.90%  '
$ws.Range("D41").Value = '1.98'
$ws.Range("E41").Value = '  -6.08%  '
$ws.Range("E42").Value = '  -5.24%  '
$ws.Range("D43").Value = '0.291'
$ws.Range("E43").Value = '  -5.63%  '
$ws.Range("D44").Value = '40.47'
$ws.Range("E44").Value = '  -8.52%  '
$ws.Range("D45").Value = '0.0348'
$ws.Range("E45").Value = '  -3.27%  '
$ws.Range("D46").Value = '371.60'
$ws.Range("E46").Value = '  -5.54%  '
$ws.Range("D47").Value = '2.695.92'
$ws.Range("E47").Value = '  -0.84%  '
$ws.Range("D48").Value = '132.14'
$ws.Range("E48").Value = '  -1.07%  '
$ws.Range("D50").Value = '23.94'
$ws.Range("E50").Value = '  -10.34%  '
$ws.Range("E51").Value = '  -2.52%  '

foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
